# Reproduce the edit: the user selected A3:A7 (the part-number values below
# the header/first row) and deleted their contents, leaving the cell
# formatting (style) intact. The selection is left spanning A3:A7 afterward.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$rng = $ws.Range("A3:A7")
[void]$rng.ClearContents()
[void]$rng.Select()
